# Payment with Card Test Cases completed
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Payment with Card" test case row (row 4): card number + expiry,
# entered as text (leading apostrophe forces text / quotePrefix, same as
# typing '4111111111111111 and '0222 directly into the cells).
$ws.Range("A4").Value = "'4111111111111111"
$ws.Range("B4").Value = "'0222"

# Column A now needs to be wide enough to fit the card number, same as the
# other "best fit" columns on this sheet.
$ws.Columns.Item(1).AutoFit() | Out-Null

# Excel leaves the newly-entered row selected (whole row) after data entry.
$ws.Rows.Item(4).Select() | Out-Null

# Restore/update the workbook window position recorded in the last view.
$win = $wb.Windows.Item(1)
$win.Left = 1300
$win.Top = 3680
